# Auto-generated PowerShell (iron_native COM-interop) script
# Adds two new handback entries (rows 38/39) to Overview, zh-cn and de-de sheets

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# ---- Entity 644510b7-72de-41c5-b953-cfb0bd0c6023 (row 38) ----
$wsOverview.Hyperlinks.Add($wsOverview.Range("A38"), "https://github.com/OpenLocalizationTest/oltest/blob/13c4f007f56ff608b29928e42ebb5e8d2246a2e6/e2e/644510b7-72de-41c5-b953-cfb0bd0c6023.md", "", "", "644510b7-72de-41c5-b953-cfb0bd0c6023.md") | Out-Null
$wsOverview.Range("B38").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C38").Value = "Handed back: in sync with en-US"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A38"), "https://github.com/OpenLocalizationTest/oltest/blob/13c4f007f56ff608b29928e42ebb5e8d2246a2e6/e2e/644510b7-72de-41c5-b953-cfb0bd0c6023.md", "", "", "644510b7-72de-41c5-b953-cfb0bd0c6023.md") | Out-Null
$wsZhCn.Range("B38").Value = "Handed back: in sync with en-US"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C38"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/13c4f007f56ff608b29928e42ebb5e8d2246a2e6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/644510b7-72de-41c5-b953-cfb0bd0c6023.13c4f007f56ff608b29928e42ebb5e8d2246a2e6.zh-cn.xlf", "", "", "644510b7-72de-41c5-b953-cfb0bd0c6023.13c4f007f56ff608b29928e42ebb5e8d2246a2e6.zh-cn.xlf") | Out-Null
$wsZhCn.Range("D38").Value = "2016-03-03 13:15:12"
$wsZhCn.Range("D38").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E38"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/13c4f007f56ff608b29928e42ebb5e8d2246a2e6/e2e/644510b7-72de-41c5-b953-cfb0bd0c6023.md", "", "", "644510b7-72de-41c5-b953-cfb0bd0c6023.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F38"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/13c4f007f56ff608b29928e42ebb5e8d2246a2e6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/644510b7-72de-41c5-b953-cfb0bd0c6023.13c4f007f56ff608b29928e42ebb5e8d2246a2e6.zh-cn.xlf", "", "", "644510b7-72de-41c5-b953-cfb0bd0c6023.13c4f007f56ff608b29928e42ebb5e8d2246a2e6.zh-cn.xlf") | Out-Null
$wsZhCn.Range("G38").Value = "2016-03-03 13:16:17"
$wsZhCn.Range("H38").Value = "Include"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A38"), "https://github.com/OpenLocalizationTest/oltest/blob/13c4f007f56ff608b29928e42ebb5e8d2246a2e6/e2e/644510b7-72de-41c5-b953-cfb0bd0c6023.md", "", "", "644510b7-72de-41c5-b953-cfb0bd0c6023.md") | Out-Null
$wsDeDe.Range("B38").Value = "Handed back: in sync with en-US"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C38"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/13c4f007f56ff608b29928e42ebb5e8d2246a2e6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/644510b7-72de-41c5-b953-cfb0bd0c6023.13c4f007f56ff608b29928e42ebb5e8d2246a2e6.de-de.xlf", "", "", "644510b7-72de-41c5-b953-cfb0bd0c6023.13c4f007f56ff608b29928e42ebb5e8d2246a2e6.de-de.xlf") | Out-Null
$wsDeDe.Range("D38").Value = "2016-03-03 13:15:24"
$wsDeDe.Range("D38").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E38"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/13c4f007f56ff608b29928e42ebb5e8d2246a2e6/e2e/644510b7-72de-41c5-b953-cfb0bd0c6023.md", "", "", "644510b7-72de-41c5-b953-cfb0bd0c6023.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F38"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/13c4f007f56ff608b29928e42ebb5e8d2246a2e6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/644510b7-72de-41c5-b953-cfb0bd0c6023.13c4f007f56ff608b29928e42ebb5e8d2246a2e6.de-de.xlf", "", "", "644510b7-72de-41c5-b953-cfb0bd0c6023.13c4f007f56ff608b29928e42ebb5e8d2246a2e6.de-de.xlf") | Out-Null
$wsDeDe.Range("G38").Value = "2016-03-03 13:16:45"
$wsDeDe.Range("H38").Value = "Include"

# ---- Entity 98644877-9fba-49f4-988d-f79b54b3127b (row 39) ----
$wsOverview.Hyperlinks.Add($wsOverview.Range("A39"), "https://github.com/OpenLocalizationTest/oltest/blob/abe0e60c061ec476a0ed333ccd0a434083174562/e2e/98644877-9fba-49f4-988d-f79b54b3127b.md", "", "", "98644877-9fba-49f4-988d-f79b54b3127b.md") | Out-Null
$wsOverview.Range("B39").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C39").Value = "Handed back: in sync with en-US"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A39"), "https://github.com/OpenLocalizationTest/oltest/blob/abe0e60c061ec476a0ed333ccd0a434083174562/e2e/98644877-9fba-49f4-988d-f79b54b3127b.md", "", "", "98644877-9fba-49f4-988d-f79b54b3127b.md") | Out-Null
$wsZhCn.Range("B39").Value = "Handed back: in sync with en-US"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C39"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/abe0e60c061ec476a0ed333ccd0a434083174562/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/98644877-9fba-49f4-988d-f79b54b3127b.abe0e60c061ec476a0ed333ccd0a434083174562.zh-cn.xlf", "", "", "98644877-9fba-49f4-988d-f79b54b3127b.abe0e60c061ec476a0ed333ccd0a434083174562.zh-cn.xlf") | Out-Null
$wsZhCn.Range("D39").Value = "2016-03-03 13:15:12"
$wsZhCn.Range("D39").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E39"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/abe0e60c061ec476a0ed333ccd0a434083174562/e2e/98644877-9fba-49f4-988d-f79b54b3127b.md", "", "", "98644877-9fba-49f4-988d-f79b54b3127b.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F39"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/abe0e60c061ec476a0ed333ccd0a434083174562/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/98644877-9fba-49f4-988d-f79b54b3127b.abe0e60c061ec476a0ed333ccd0a434083174562.zh-cn.xlf", "", "", "98644877-9fba-49f4-988d-f79b54b3127b.abe0e60c061ec476a0ed333ccd0a434083174562.zh-cn.xlf") | Out-Null
$wsZhCn.Range("G39").Value = "2016-03-03 13:16:17"
$wsZhCn.Range("H39").Value = "Include"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A39"), "https://github.com/OpenLocalizationTest/oltest/blob/abe0e60c061ec476a0ed333ccd0a434083174562/e2e/98644877-9fba-49f4-988d-f79b54b3127b.md", "", "", "98644877-9fba-49f4-988d-f79b54b3127b.md") | Out-Null
$wsDeDe.Range("B39").Value = "Handed back: in sync with en-US"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C39"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/abe0e60c061ec476a0ed333ccd0a434083174562/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/98644877-9fba-49f4-988d-f79b54b3127b.abe0e60c061ec476a0ed333ccd0a434083174562.de-de.xlf", "", "", "98644877-9fba-49f4-988d-f79b54b3127b.abe0e60c061ec476a0ed333ccd0a434083174562.de-de.xlf") | Out-Null
$wsDeDe.Range("D39").Value = "2016-03-03 13:15:24"
$wsDeDe.Range("D39").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E39"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/abe0e60c061ec476a0ed333ccd0a434083174562/e2e/98644877-9fba-49f4-988d-f79b54b3127b.md", "", "", "98644877-9fba-49f4-988d-f79b54b3127b.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F39"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/abe0e60c061ec476a0ed333ccd0a434083174562/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/98644877-9fba-49f4-988d-f79b54b3127b.abe0e60c061ec476a0ed333ccd0a434083174562.de-de.xlf", "", "", "98644877-9fba-49f4-988d-f79b54b3127b.abe0e60c061ec476a0ed333ccd0a434083174562.de-de.xlf") | Out-Null
$wsDeDe.Range("G39").Value = "2016-03-03 13:16:45"
$wsDeDe.Range("H39").Value = "Include"
